$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 40/41: PaxDollar and RenderToken swap positions, with updated data
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D40" "1.915"
$ws.Range("E40").Value = "  -1.85%  "

$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D41" "1.002"
$ws.Range("E41").Value = "  +0.24%  "

# Price (D) and Volume(1h) (E) updates for remaining rows
Set-TextValue "D2" "25.855.45"
$ws.Range("E2").Value = "  +0.06%  "
Set-TextValue "D3" "1.737.14"
$ws.Range("E3").Value = "  +0.18%  "
Set-TextValue "D4" "1.001"
$ws.Range("E4").Value = "  +0.10%  "
Set-TextValue "D5" "237.87"
$ws.Range("E5").Value = "  +2.86%  "
Set-TextValue "D6" "1.002"
$ws.Range("E6").Value = "  +0.18%  "
Set-TextValue "D7" "0.5140"
$ws.Range("E7").Value = "  -0.09%  "
Set-TextValue "D8" "0.2717"
$ws.Range("E8").Value = "  -2.35%  "
Set-TextValue "D9" "39.70"
$ws.Range("E9").Value = "  +0.90%  "
Set-TextValue "D10" "0.06105"
$ws.Range("E10").Value = "  -0.09%  "
Set-TextValue "D11" "1.739.38"
$ws.Range("E11").Value = "  +0.26%  "
Set-TextValue "D12" "0.07179"
$ws.Range("E12").Value = "  +2.22%  "
$ws.Range("E13").Value = "  -2.07%  "
Set-TextValue "D14" "0.6354"
$ws.Range("E14").Value = "  -1.08%  "
Set-TextValue "D15" "4.587"
$ws.Range("E15").Value = "  +1.28%  "
Set-TextValue "D16" "77.11"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("E18").Value = "  +0.12%  "
Set-TextValue "D19" "25.867.43"
$ws.Range("E19").Value = "  +0.14%  "
Set-TextValue "D20" "11.70"
$ws.Range("E20").Value = "  +1.93%  "
Set-TextValue "D21" "0.000006695"
$ws.Range("E21").Value = "  +1.14%  "
Set-TextValue "D22" "1.959.59"
$ws.Range("E22").Value = "  -0.51%  "
Set-TextValue "D23" "4.278"
$ws.Range("E23").Value = "  +3.21%  "
Set-TextValue "D24" "8.630"
$ws.Range("E24").Value = "  -1.27%  "
Set-TextValue "D25" "5.229"
$ws.Range("E25").Value = "  +2.07%  "
Set-TextValue "D26" "139.22"
$ws.Range("E26").Value = "  -0.29%  "
Set-TextValue "D27" "1.515"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("E28").Value = "  +0.93%  "
Set-TextValue "D29" "1.761"
$ws.Range("E29").Value = "  -1.84%  "
Set-TextValue "D30" "105.58"
$ws.Range("E30").Value = "  +3.52%  "
Set-TextValue "D31" "3.883"
$ws.Range("E31").Value = "  +5.20%  "
Set-TextValue "D32" "0.08338"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("E33").Value = "  +4.47%  "
$ws.Range("E34").Value = "  +1.72%  "
Set-TextValue "D35" "2.634"
$ws.Range("E35").Value = "  +0.66%  "
Set-TextValue "D36" "0.9830"
$ws.Range("E36").Value = "  +0.33%  "
Set-TextValue "D37" "0.6203"
$ws.Range("E37").Value = "  +1.31%  "
Set-TextValue "D38" "2.698"
$ws.Range("E38").Value = "  +2.26%  "
Set-TextValue "D39" "0.01592"
$ws.Range("E39").Value = "  +0.88%  "
Set-TextValue "D42" "97.68"
$ws.Range("E42").Value = "  -2.72%  "
$ws.Range("E43").Value = "  +0.56%  "
Set-TextValue "D44" "0.7345"
$ws.Range("E44").Value = "  +1.21%  "
Set-TextValue "D45" "4.932"
$ws.Range("E45").Value = "  -0.56%  "
Set-TextValue "D46" "0.1132"
$ws.Range("E46").Value = "  +0.86%  "
Set-TextValue "D47" "0.05280"
$ws.Range("E47").Value = "  -2.08%  "
Set-TextValue "D48" "6.181"
$ws.Range("E48").Value = "  -1.22%  "
Set-TextValue "D49" "54.57"
$ws.Range("E49").Value = "  +3.14%  "
Set-TextValue "D50" "30.49"
$ws.Range("E50").Value = "  +1.41%  "
Set-TextValue "D51" "7.539"
$ws.Range("E51").Value = "  -0.53%  "
